$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: "center of mass adjustment" now subtracts the b1 (B3/C3) position,
# fixing the center-of-mass calc to exclude/offset boid b1 itself.
$ws.Range("B22").Formula = "=(B2+B4+B5)/3-B3"
$ws.Range("C22").Formula = "=(C2+C4+C5)/3-C3"

# Row 30 / E30 weight tweak (avoidance-adjustment weight)
$ws.Range("E30").Value = 0.1

# Row 37: center of mass adjustment formulas for the "For B5" section now
# offset by the bound values (-6.5 / 6) before scaling.
$ws.Range("C37").Formula = "=0.01*((SUM(B2:B5)/4)-(-6.5))"
$ws.Range("D37").Formula = "=0.01*((SUM(C2:C5)/4)-6)"

# Update the active cell selection to match the saved view state.
$ws.Range("D38").Select()
